# Updated symbol list on Mon Jan  2 22:26:48 UTC 2023 with GitHub Actions
# Refreshes price/volume figures and re-syncs a few reordered coin rows
# (rows 9-14 and 41-43) to match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.01"
$ws.Range("E2").Value = "'0.90%"
$ws.Range("D3").Value = "'29.50"
$ws.Range("E3").Value = "'7.35%"
$ws.Range("D4").Value = "'5.191"
$ws.Range("E4").Value = "'1.49%"
$ws.Range("D5").Value = "'0.05732"
$ws.Range("E5").Value = "'0.83%"
$ws.Range("D6").Value = "'6.561"
$ws.Range("E6").Value = "'0.61%"
$ws.Range("D7").Value = "'0.8579"
$ws.Range("E7").Value = "'4.61%"
$ws.Range("D8").Value = "'0.8671"
$ws.Range("E8").Value = "'1.86%"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'0.1366"
$ws.Range("E9").Value = "'2.44%"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = "'0.07084"
$ws.Range("E10").Value = "'2.09%"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = "'0.02996"
$ws.Range("E11").Value = "'3.94%"
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = "'0.09385"
$ws.Range("E12").Value = "'-0.07%"
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = "'0.001531"
$ws.Range("E13").Value = "'0.15%"
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = "'0.0006000"
$ws.Range("E14").Value = "'-0.34%"
$ws.Range("D15").Value = "'0.006076"
$ws.Range("E15").Value = "'-2.24%"
$ws.Range("E16").Value = "'5,224.37%"
$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'-0.51%"
$ws.Range("D18").Value = "'3.101"
$ws.Range("E18").Value = "'2.94%"
$ws.Range("D19").Value = "'2.187"
$ws.Range("E19").Value = "'-1.94%"
$ws.Range("E20").Value = "'0.53%"
$ws.Range("D21").Value = "'0.03306"
$ws.Range("E21").Value = "'3.04%"
$ws.Range("E22").Value = "'1.24%"
$ws.Range("D23").Value = "'3.469"
$ws.Range("E23").Value = "'-2.60%"
$ws.Range("D24").Value = "'0.04147"
$ws.Range("E24").Value = "'3.15%"
$ws.Range("E25").Value = "'0.43%"
$ws.Range("E26").Value = "'0.88%"
$ws.Range("D27").Value = "'0.004994"
$ws.Range("E27").Value = "'11.51%"
$ws.Range("E28").Value = "'2.52%"
$ws.Range("E40").Value = "'1.09%"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = "'0.1072"
$ws.Range("E41").Value = "'1.16%"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.002449"
$ws.Range("E42").Value = "'4.20%"
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.003519"
$ws.Range("E43").Value = "'-40.45%"
$ws.Range("D44").Value = "'0.009486"
$ws.Range("D45").Value = "'0.00005264"
$ws.Range("E45").Value = "'3.19%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D47").Value = "'0.05700"
$ws.Range("E47").Value = "'-43.57%"
$ws.Range("D48").Value = "'0.002278"
$ws.Range("E48").Value = "'-9.42%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.00%"
